# Update the timestamp column (A2:A11) on the "ランサーズ" sheet to reflect
# the new scrape/acquisition time, per commit: "Append: 2025-09-28 01:50 JST"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-28 01:50:26"

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
